# Weekly update: insert a new price record as row 534, pushing the
# existing rows 534:599 down to 535:600.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("534:534").Insert()

$ws.Range("A534").Value = 9
$ws.Range("B534").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C534").Value = "Metropolitana"
$ws.Range("D534").Value = 45142
$ws.Range("E534").Value = 13
$ws.Range("F534").Value = 100112044
$ws.Range("G534").Value = "Perejil"
$ws.Range("H534").Value = "Sin especificar"
$ws.Range("I534").Value = "Primera"
$ws.Range("J534").Value = 70
$ws.Range("K534").Value = 13000
$ws.Range("L534").Value = 14000
$ws.Range("M534").Value = 13500
$ws.Range("N534").Value = "$/docena de atados"
$ws.Range("O534").Value = "Región Metropolitana"
$ws.Range("P534").Value = 4500
$ws.Range("Q534").Value = 3
$ws.Range("R534").Value = "Hortaliza"
